# Apply the "Top level menu, Restriction based data filtration, reports,
# new home page view and, issue fixes" edit to the cms_sample workbook.
#
# Summary of the change (per the OOXML diff):
#   1. C2 (Interviewer) gains two more co-interviewers: the cell's text
#      changes from "Little Thunder, Julie Pearson" to
#      "Little Thunder, Julie Pearson;Finchum, Tanya; Bishop, Alex"
#      (this introduces a brand-new shared string).
#   2. The sheet's view no longer scrolls to keep AA1 pinned as the
#      top-left cell, and the live selection moves from AF29 to C3 -
#      i.e. the saved view now shows the new home-page-ish layout
#      starting at column A with C3 selected.
#   3. Explicit column widths are now defined for columns B and C so the
#      (now more visible) Interviewee / Interviewer columns read well.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# 1) Update the Interviewer cell's text (this allocates the new shared
#    string automatically, just like typing the new value into Excel).
$ws.Range("C2").Value = "Little Thunder, Julie Pearson;Finchum, Tanya; Bishop, Alex"

# 2) Move the visible selection to C3 and let the window scroll back so
#    column A is the top-left cell again (drops the old topLeftCell="AA1").
$ws.Range("C3").Select()

# 3) Give columns B (Interviewee) and C (Interviewer) explicit widths so
#    the newly-lengthened interviewer text is readable.
$ws.Columns.Item(2).ColumnWidth = 24.140625
$ws.Columns.Item(3).ColumnWidth = 82.7109375
